$d = $word.ActiveDocument
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# Edit 1: split the sentence
#   "The following commands are to be run when launching the application:"
# into three runs:
#   "The following commands are to be run when launching the application"
#   ", when inside the repo’s root directory"
#   ":"
# (all three keep the same Times New Roman / kern 0 / szCs 24 formatting that
#  the original run already had).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "The following commands are to be run when launching the application:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The following commands are to be run when launching the application",
    2) | Out-Null

$launchPara = $d.Paragraphs.Item(10)
$launchEnd = $launchPara.Range
$insertPoint1 = $d.Range($launchEnd.End - 1, $launchEnd.End - 1)
$insertPoint1.InsertAfter(", when inside the repo" + $apos + "s root directory")

$launchPara2 = $d.Paragraphs.Item(10)
$launchEnd2 = $launchPara2.Range
$insertPoint2 = $d.Range($launchEnd2.End - 1, $launchEnd2.End - 1)
$insertPoint2.InsertAfter(":")

# ---------------------------------------------------------------------------
# Edit 2: after the (empty) paragraph that follows "where <port> ..." and
# precedes the trailing empty paragraph, add the new "reset the database"
# instructions plus a couple of spacer paragraphs and one more trailing
# empty paragraph.
# ---------------------------------------------------------------------------

# Locate the empty paragraph right before the final empty paragraph / sectPr.
$count = $d.Paragraphs.Count
$spacerPara = $d.Paragraphs.Item($count - 1)

# Give that spacer paragraph's mark the Times New Roman / kern 0 / szCs 24
# run-properties (matches the rest of the document's body-text paragraphs).
$spacerRange = $spacerPara.Range
$spacerRange.InsertAfter("x")
$spacerPara2 = $d.Paragraphs.Item($count - 1)
$tmpRange = $spacerPara2.Range
$tmpRange.Font.Name = "Times New Roman"
$tmpRange.Font.Kerning = 0
$tmpRange.Font.SizeBi = 12
$spacerPara3 = $d.Paragraphs.Item($count - 1)
$charRange = $d.Range($spacerPara3.Range.Start, $spacerPara3.Range.Start + 1)
$charRange.Delete()

# Paragraph: "If you are testing the database and you would like to reset its
# contents, simply run "
$spacerPara4 = $d.Paragraphs.Item($count - 1)
$spacerPara4.Range.InsertParagraphAfter()
$resetIntroPara = $d.Paragraphs.Item($count)
$resetIntroRange = $resetIntroPara.Range
$resetIntroRange.Text = "If you are testing the database and you would like to reset its contents, simply run "
$resetIntroRange2 = $d.Paragraphs.Item($count).Range
$resetIntroRange2.Font.Name = "Times New Roman"
$resetIntroRange2.Font.Kerning = 0
$resetIntroRange2.Font.SizeBi = 12

# Blank paragraph (no special run formatting).
$resetIntroPara2 = $d.Paragraphs.Item($count)
$resetIntroPara2.Range.InsertParagraphAfter()

# Paragraph with the "./resetdb.sh" command (indented, Consolas font).
$blankPara = $d.Paragraphs.Item($count + 1)
$blankPara.Range.InsertParagraphAfter()
$cmdPara = $d.Paragraphs.Item($count + 2)
$cmdRange = $cmdPara.Range
$cmdRange.Text = "./resetdb.sh"
$cmdPara2 = $d.Paragraphs.Item($count + 2)
$cmdPara2.Format.LeftIndent = $word.InchesToPoints(0.33)
$cmdRange2 = $d.Paragraphs.Item($count + 2).Range
$cmdRange2.Font.Name = "Consolas"
$cmdRange2.Font.NameAscii = "Consolas"
$cmdRange2.LanguageID = 1076

# Blank paragraph after the command.
$cmdPara3 = $d.Paragraphs.Item($count + 2)
$cmdPara3.Range.InsertParagraphAfter()

# One more, fully empty trailing paragraph (no pPr contents at all).
$blankPara2 = $d.Paragraphs.Item($count + 3)
$blankPara2.Range.InsertParagraphAfter()
